# Add a new "UK" Test Data sheet, cloned from the "Poland" sheet (same
# template: column widths, styles, merged cells, page setup) and extended
# with the two extra repeater rows ("P32AR"/"P32DR") that the UK sheet
# needs but Poland doesn't, then filled in with UK-specific market data.

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Duplicate "Poland" immediately after itself - this carries over all
# formatting (column widths, cell styles, merged cells, page setup) and
# becomes the new active sheet (so it naturally gains tabSelected, and the
# previously-active sheet naturally loses it).
$poland.Copy($null, $poland) | Out-Null
$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# The UK sheet lists two more repeater models than Poland's, so insert two
# rows above the tail of the list (currently rows 16 "PR1DS"/17 "PR8AS"),
# inheriting the formatting of the row above them.
$uk.Rows.Item(16).Resize(2).Insert() | Out-Null
$uk.Range("A15").Copy() | Out-Null
$uk.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$uk.Range("A16").Value = "P32AR"
$uk.Range("A17").Value = "P32DR"

# Fill in the UK-specific values. B4 is set before B2 so the new shared
# strings are appended in that order ("NGC-2741/T3342" then "UK Market").
$uk.Range("B4").Value = "NGC-2741/T3342"
$uk.Range("B2").Value = "UK Market"

# Match the recorded view state of the new sheet (selection on B4).
$uk.Range("B4").Select() | Out-Null
